$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    ,@("2023-10-07 20:10:44", "hatespeech", "setting8", 10, 5.291274608099837)
    ,@("2023-10-07 20:10:44", "hatespeech", "setting1", 10, 5.695187342672982)
    ,@("2023-10-07 20:10:48", "hatespeech", "setting3", 10, 5.811032671043269)
    ,@("2023-10-07 20:10:50", "hatespeech", "setting5", 10, 5.542625778243734)
    ,@("2023-10-07 20:10:54", "hatespeech", "setting2", 10, 5.677348860833649)
    ,@("2023-10-07 20:10:54", "hatespeech", "setting7", 10, 5.416487599006271)
    ,@("2023-10-07 20:10:55", "hatespeech", "def", 10, 5.498188539206063)
    ,@("2023-10-07 20:10:55", "hatespeech", "setting6", 20, 5.466001344594968)
    ,@("2023-10-07 20:10:55", "hatespeech", "setting5", 20, 5.431023151521125)
    ,@("2023-10-07 20:10:55", "hatespeech", "setting4", 10, 5.376840756031892)
    ,@("2023-10-07 20:10:55", "hatespeech", "setting3", 40, 5.676024509278975)
    ,@("2023-10-07 20:10:55", "hatespeech", "setting1", 40, 5.432202026289972)
    ,@("2023-10-07 20:10:56", "hatespeech", "def", 30, 5.261309785235558)
    ,@("2023-10-07 20:10:56", "hatespeech", "setting4", 40, 5.167869402747034)
    ,@("2023-10-07 20:10:56", "hatespeech", "def", 40, 5.318752964990559)
    ,@("2023-10-07 20:10:56", "hatespeech", "setting3", 30, 5.737660878335407)
    ,@("2023-10-07 20:10:56", "hatespeech", "setting2", 20, 5.58983429363817)
    ,@("2023-10-07 20:10:56", "hatespeech", "setting3", 20, 5.713483319319765)
    ,@("2023-10-07 20:10:57", "hatespeech", "setting2", 40, 5.537689107541402)
    ,@("2023-10-07 20:10:57", "hatespeech", "setting1", 30, 5.490376847492876)
    ,@("2023-10-07 20:10:57", "hatespeech", "setting4", 30, 5.228675362844228)
    ,@("2023-10-07 20:10:57", "hatespeech", "setting7", 30, 5.293886123450595)
    ,@("2023-10-07 20:10:57", "hatespeech", "setting8", 20, 5.213577271181734)
    ,@("2023-10-07 20:10:57", "hatespeech", "setting1", 20, 5.467478110137145)
    ,@("2023-10-07 20:10:57", "hatespeech", "setting6", 10, 5.410521191842305)
    ,@("2023-10-07 20:10:57", "hatespeech", "def", 20, 5.311152804181749)
    ,@("2023-10-07 20:10:58", "hatespeech", "setting6", 30, 5.326747805734123)
    ,@("2023-10-07 20:10:58", "hatespeech", "setting5", 30, 5.403353644352125)
    ,@("2023-10-07 20:10:58", "hatespeech", "setting8", 30, 5.172802390874662)
    ,@("2023-10-07 20:10:58", "hatespeech", "setting4", 20, 5.401152790140427)
    ,@("2023-10-07 20:10:58", "hatespeech", "setting7", 20, 5.467592929978188)
    ,@("2023-10-07 20:10:58", "hatespeech", "setting2", 30, 5.474697923367764)
    ,@("2023-10-07 20:11:14", "hatespeech", "setting5", 40, 5.428426457365942)
    ,@("2023-10-07 20:11:14", "hatespeech", "setting6", 40, 5.442066440598762)
    ,@("2023-10-07 20:11:16", "hatespeech", "setting7", 40, 5.359746753566184)
    ,@("2023-10-07 20:11:16", "hatespeech", "setting8", 40, 5.088884102083386)
)

$startRow = 256
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
